$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.959.06'
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").Value = '3.213.43'
$ws.Range("E3").Value = '  +1.33%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.212.48'
$ws.Range("E8").Value = '  +1.30%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("E10").Value = '  -1.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.38%  '

$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.31%  '

$ws.Range("D15").Value = '3.739.94'
$ws.Range("E15").Value = '  +1.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.47'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = '66.114.84'
$ws.Range("E17").Value = '  +1.73%  '

$ws.Range("D18").Value = '3.213.20'
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '510.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.735'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.52%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +2.96%  '

$ws.Range("E28").Value = '  +1.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.35%  '

$ws.Range("E30").Value = '  +3.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.27%  '

$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("E34").Value = '  +0.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.61'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0905'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.297'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.17%  '

$ws.Range("E43").Value = '  +0.60%  '

$ws.Range("D44").Value = '2.954.06'
$ws.Range("E44").Value = '  -3.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.38%  '

$ws.Range("E46").Value = '  +5.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.76%  '

$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
